$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so values like "1.00" or "25.10"
# are not coerced into numbers and lose their trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('B48').NumberFormat = '@'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('B49').NumberFormat = '@'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '67.668.65'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '2.477.67'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '585.78'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '173.39'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.512'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').Value = '0.143'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '2.935.55'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').Value = '25.10'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '67.601.00'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '2.478.32'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '10.74'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').Value = '7.33'
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('D20').Value = '345.31'
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = '4.08'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '70.62'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').Value = '4.15'
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('D25').Value = '1.66'
$ws.Range('E25').Value = '  -8.62%  '
$ws.Range('D26').Value = '8.78'
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '0.0₃0881'
$ws.Range('E29').Value = '  -3.35%  '
$ws.Range('D30').Value = '494.16'
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('D31').Value = '7.66'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').Value = '1.23'
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').Value = '164.49'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('D36').Value = '0.119'
$ws.Range('E36').Value = '  +0.96%  '
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').Value = '18.20'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  -3.33%  '
$ws.Range('D41').Value = '1.70'
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('D42').Value = '0.321'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('D43').Value = '4.72'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('D44').Value = '2.35'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').Value = '147.03'
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('D46').Value = '3.49'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D47').Value = '0.507'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0734'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0250'
$ws.Range('E49').Value = '  -5.47%  '
$ws.Range('D50').Value = '1.54'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').Value = '0.574'
$ws.Range('E51').Value = '  -1.62%  '
